# Update market price columns (H-N) across multiple sheets with refreshed
# Universalis price data, as produced by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 557077.4399999999
$ws.Range("I19").Value = 2000160.2
$ws.Range("J19").Value = 2045.6154
$ws.Range("K19").Value = 2000160.2
$ws.Range("L19").Value = 2045.6154
$ws.Range("M19").Value = -1999985.2
$ws.Range("N19").Value = -2395.6154
$ws.Range("H76").Value = 3450
$ws.Range("J76").Value = 3450
$ws.Range("L76").Value = 3450
$ws.Range("N76").Value = -4080
$ws.Range("H79").Value = 3450
$ws.Range("J79").Value = 3450
$ws.Range("L79").Value = 3450
$ws.Range("N79").Value = -5634
$ws.Range("I132").Value = 1142.619
$ws.Range("J132").Value = 1478
$ws.Range("K132").Value = 3427.857
$ws.Range("L132").Value = 4434
$ws.Range("M132").Value = -897.857
$ws.Range("N132").Value = -9494
$ws.Range("H138").Value = 2558.2754
$ws.Range("J138").Value = 2074.5918
$ws.Range("L138").Value = 6223.7754
$ws.Range("N138").Value = -16503.7754
$ws.Range("H139").Value = 71797.5
$ws.Range("J139").Value = 71797.5
$ws.Range("L139").Value = 71797.5
$ws.Range("N139").Value = -82077.5
$ws.Range("H140").Value = 81251.78999999999
$ws.Range("J140").Value = 81251.78999999999
$ws.Range("L140").Value = 81251.78999999999
$ws.Range("N140").Value = -91611.78999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4229.7
$ws.Range("I32").Value = 3369.691
$ws.Range("K32").Value = 3369.691
$ws.Range("M32").Value = -3082.691
$ws.Range("H45").Value = 1579.4667
$ws.Range("I45").Value = 1050
$ws.Range("J45").Value = 1772
$ws.Range("K45").Value = 1050
$ws.Range("L45").Value = 1772
$ws.Range("M45").Value = -673
$ws.Range("N45").Value = -2526
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("H97").Value = 1305
$ws.Range("I97").Value = 1217.4445
$ws.Range("K97").Value = 1217.4445
$ws.Range("M97").Value = -721.4445000000001
$ws.Range("H109").Value = 61962.332
$ws.Range("J109").Value = 61962.332
$ws.Range("L109").Value = 61962.332
$ws.Range("N109").Value = -64736.332
$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1044.375
$ws.Range("K122").Value = 3133.125
$ws.Range("M122").Value = -683.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 79599.5
$ws.Range("I86").Value = 3137.6667
$ws.Range("K86").Value = 3137.6667
$ws.Range("M86").Value = -2014.6667
$ws.Range("H89").Value = 79599.5
$ws.Range("I89").Value = 3137.6667
$ws.Range("K89").Value = 15688.3335
$ws.Range("M89").Value = -10072.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2470.7273
$ws.Range("I31").Value = 2041
$ws.Range("K31").Value = 2041
$ws.Range("M31").Value = -1746
$ws.Range("H34").Value = 2470.7273
$ws.Range("I34").Value = 2041
$ws.Range("K34").Value = 2041
$ws.Range("M34").Value = -1839
$ws.Range("H105").Value = 982.3
$ws.Range("I105").Value = 988.375
$ws.Range("K105").Value = 988.375
$ws.Range("M105").Value = 758.625
$ws.Range("H122").Value = 4750.231
$ws.Range("I122").Value = 3193.5557
$ws.Range("K122").Value = 9580.667099999999
$ws.Range("M122").Value = -7130.667099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 156.76923
$ws.Range("J12").Value = 219.375
$ws.Range("L12").Value = 658.125
$ws.Range("N12").Value = -1004.125
$ws.Range("H113").Value = 5291.364
$ws.Range("J113").Value = 746.7222
$ws.Range("L113").Value = 2240.1666
$ws.Range("N113").Value = -6580.1666
$ws.Range("H131").Value = 17043.139
$ws.Range("J131").Value = 19198.053
$ws.Range("L131").Value = 57594.159
$ws.Range("N131").Value = -67674.159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13874.875
$ws.Range("I70").Value = 29599.666
$ws.Range("K70").Value = 29599.666
$ws.Range("M70").Value = -29329.666
$ws.Range("H73").Value = 13874.875
$ws.Range("I73").Value = 29599.666
$ws.Range("K73").Value = 29599.666
$ws.Range("M73").Value = -28663.666
$ws.Range("H102").Value = 2215.1667
$ws.Range("I102").Value = 1847.3529
$ws.Range("K102").Value = 1847.3529
$ws.Range("M102").Value = -225.3529000000001
$ws.Range("H113").Value = 859.9231
$ws.Range("I113").Value = 393.57144
$ws.Range("K113").Value = 393.57144
$ws.Range("M113").Value = 1776.42856
$ws.Range("H122").Value = 1056.8182
$ws.Range("I122").Value = 862
$ws.Range("J122").Value = 1290.6
$ws.Range("K122").Value = 2586
$ws.Range("L122").Value = 3871.8
$ws.Range("M122").Value = -136
$ws.Range("N122").Value = -8771.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1544.0834
$ws.Range("I22").Value = 1157.8334
$ws.Range("J22").Value = 1930.3334
$ws.Range("K22").Value = 1157.8334
$ws.Range("L22").Value = 1930.3334
$ws.Range("M22").Value = -862.8334
$ws.Range("N22").Value = -2520.3334
$ws.Range("H27").Value = 1544.0834
$ws.Range("I27").Value = 1157.8334
$ws.Range("J27").Value = 1930.3334
$ws.Range("K27").Value = 1157.8334
$ws.Range("L27").Value = 1930.3334
$ws.Range("M27").Value = -1050.8334
$ws.Range("N27").Value = -2144.3334
$ws.Range("H40").Value = 10911.444
$ws.Range("I40").Value = 11000.538
$ws.Range("J40").Value = 10679.8
$ws.Range("K40").Value = 11000.538
$ws.Range("L40").Value = 10679.8
$ws.Range("M40").Value = -10864.538
$ws.Range("N40").Value = -10951.8
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null
$ws.Range("H93").Value = 16667175
$ws.Range("J93").Value = 66667216
$ws.Range("L93").Value = 66667216
$ws.Range("N93").Value = -66669712
$ws.Range("H132").Value = 3237
$ws.Range("I132").Value = 1963.125
$ws.Range("J132").Value = 4020.923
$ws.Range("K132").Value = 5889.375
$ws.Range("L132").Value = 12062.769
$ws.Range("M132").Value = -3359.375
$ws.Range("N132").Value = -17122.769
$ws.Range("H136").Value = 2420.8
$ws.Range("J136").Value = 2667.5
$ws.Range("L136").Value = 8002.5
$ws.Range("N136").Value = -13102.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1103.8334
$ws.Range("I100").Value = 924.8
$ws.Range("K100").Value = 1849.6
$ws.Range("M100").Value = -1308.6
$ws.Range("H122").Value = 20463.414
$ws.Range("I122").Value = 31245.115
$ws.Range("J122").Value = 1775.1333
$ws.Range("K122").Value = 93735.345
$ws.Range("L122").Value = 5325.3999
$ws.Range("M122").Value = -91285.345
$ws.Range("N122").Value = -10225.3999
